$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a numeric-looking word (e.g. "2", "3") into a cell as TEXT,
# while preserving the cells existing number/border/font formatting.
function Set-TextValue($range, $text) {
    $helper = $ws.Range("ZZ500")
    $helper.NumberFormat = "@"
    $helper.Value2 = $text
    $helper.Copy()
    $range.PasteSpecial(-4163)
    $helper.Clear()
}

# --- Update existing rows 3-21: columns A-H and J-Q (word lists re-ranked, counts updated) ---
# Row 3
$ws.Range("A3").Value = "poorly"
$ws.Range("B3").Value = 0.9782608695652174
$ws.Range("C3").Value = 45
$ws.Range("D3").Value = 45
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = $false
$ws.Range("H3").Value = 1
$ws.Range("J3").Value = "wonderful"
$ws.Range("K3").Value = 0.8928571428571429
$ws.Range("L3").Value = 50
$ws.Range("M3").Value = 50
$ws.Range("N3").Value = 1
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = $false
$ws.Range("Q3").Value = 6

# Row 4
$ws.Range("A4").Value = "disappointing"
$ws.Range("B4").Value = 0.8863636363636364
$ws.Range("C4").Value = 39
$ws.Range("D4").Value = 39
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = $false
$ws.Range("H4").Value = 5
$ws.Range("J4").Value = "awesome"
$ws.Range("K4").Value = 0.8769230769230769
$ws.Range("L4").Value = 57
$ws.Range("M4").Value = 57
$ws.Range("N4").Value = 1
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = $false
$ws.Range("Q4").Value = 8

# Row 5
$ws.Range("A5").Value = "poor"
$ws.Range("B5").Value = 0.7605633802816901
$ws.Range("C5").Value = 54
$ws.Range("D5").Value = 54
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = $false
$ws.Range("H5").Value = 17
$ws.Range("J5").Value = "excellent"
$ws.Range("K5").Value = 0.796875
$ws.Range("L5").Value = 51
$ws.Range("M5").Value = 51
$ws.Range("N5").Value = 1
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = $false
$ws.Range("Q5").Value = 13

# Row 6
$ws.Range("A6").Value = "disappointed"
$ws.Range("B6").Value = 0.7096774193548387
$ws.Range("C6").Value = 132
$ws.Range("D6").Value = 132
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = $false
$ws.Range("H6").Value = 54
$ws.Range("J6").Value = "favorite"
$ws.Range("K6").Value = 0.7956989247311828
$ws.Range("L6").Value = 74
$ws.Range("M6").Value = 74
$ws.Range("N6").Value = 1
$ws.Range("O6").Value = 0
$ws.Range("P6").Value = $false
$ws.Range("Q6").Value = 19

# Row 7
$ws.Range("A7").Value = "however"
$ws.Range("B7").Value = 0.6875
$ws.Range("C7").Value = 44
$ws.Range("D7").Value = 44
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = $false
$ws.Range("H7").Value = 20
$ws.Range("J7").Value = "classic"
$ws.Range("K7").Value = 0.6792452830188679
$ws.Range("L7").Value = 36
$ws.Range("M7").Value = 36
$ws.Range("N7").Value = 1
$ws.Range("O7").Value = 0
$ws.Range("P7").Value = $false
$ws.Range("Q7").Value = 17

# Row 8
$ws.Range("A8").Value = "waste"
$ws.Range("B8").Value = 0.6621621621621622
$ws.Range("C8").Value = 98
$ws.Range("D8").Value = 98
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = $false
$ws.Range("H8").Value = 50
$ws.Range("J8").Value = "love"
$ws.Range("K8").Value = 0.5710186513629842
$ws.Range("L8").Value = 398
$ws.Range("M8").Value = 398
$ws.Range("N8").Value = 1
$ws.Range("O8").Value = 0
$ws.Range("P8").Value = $false
$ws.Range("Q8").Value = 299

# Row 9
$ws.Range("A9").Value = "broke"
$ws.Range("B9").Value = 0.6504854368932039
$ws.Range("C9").Value = 134
$ws.Range("D9").Value = 134
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = $false
$ws.Range("H9").Value = 72
$ws.Range("J9").Value = "thank"
$ws.Range("K9").Value = 0.5507246376811594
$ws.Range("L9").Value = 38
$ws.Range("M9").Value = 38
$ws.Range("N9").Value = 1
$ws.Range("O9").Value = 0
$ws.Range("P9").Value = $false
$ws.Range("Q9").Value = 31

# Row 10
$ws.Range("A10").Value = "smaller"
$ws.Range("B10").Value = 0.5882352941176471
$ws.Range("C10").Value = 70
$ws.Range("D10").Value = 70
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = $false
$ws.Range("H10").Value = 49
$ws.Range("J10").Value = "loves"
$ws.Range("K10").Value = 0.5062240663900415
$ws.Range("L10").Value = 244
$ws.Range("M10").Value = 244
$ws.Range("N10").Value = 1
$ws.Range("O10").Value = 0
$ws.Range("P10").Value = $false
$ws.Range("Q10").Value = 238

# Row 11
$ws.Range("A11").Value = "junk"
$ws.Range("B11").Value = 0.5818181818181818
$ws.Range("C11").Value = 32
$ws.Range("D11").Value = 32
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = $false
$ws.Range("H11").Value = 23
$ws.Range("J11").Value = "great"
$ws.Range("K11").Value = 0.4581967213114754
$ws.Range("L11").Value = 559
$ws.Range("M11").Value = 559
$ws.Range("N11").Value = 1
$ws.Range("O11").Value = 0
$ws.Range("P11").Value = $false
$ws.Range("Q11").Value = 661

# Row 12
$ws.Range("A12").Value = "small"
$ws.Range("B12").Value = 0.5101449275362319
$ws.Range("C12").Value = 176
$ws.Range("D12").Value = 176
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = $false
$ws.Range("H12").Value = 169
$ws.Range("J12").Value = "loved"
$ws.Range("K12").Value = 0.363914373088685
$ws.Range("L12").Value = 119
$ws.Range("M12").Value = 119
$ws.Range("N12").Value = 1
$ws.Range("O12").Value = 0
$ws.Range("P12").Value = $false
$ws.Range("Q12").Value = 208

# Row 13
$ws.Range("A13").Value = "plastic"
$ws.Range("B13").Value = 0.3937007874015748
$ws.Range("C13").Value = 50
$ws.Range("D13").Value = 50
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = $false
$ws.Range("H13").Value = 77
$ws.Range("J13").Value = "friends"
$ws.Range("K13").Value = 0.3492063492063492
$ws.Range("L13").Value = 66
$ws.Range("M13").Value = 66
$ws.Range("N13").Value = 1
$ws.Range("O13").Value = 0
$ws.Range("P13").Value = $false
$ws.Range("Q13").Value = 123

# Row 14
$ws.Range("A14").Value = "cheap"
$ws.Range("B14").Value = 0.3886255924170616
$ws.Range("C14").Value = 82
$ws.Range("D14").Value = 82
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = $false
$ws.Range("H14").Value = 129
$ws.Range("J14").Value = "perfect"
$ws.Range("K14").Value = 0.3433734939759036
$ws.Range("L14").Value = 57
$ws.Range("M14").Value = 57
$ws.Range("N14").Value = 1
$ws.Range("O14").Value = 0
$ws.Range("P14").Value = $false
$ws.Range("Q14").Value = 109

# Row 15
$ws.Range("A15").Value = "broken"
$ws.Range("B15").Value = 0.3855421686746988
$ws.Range("C15").Value = 32
$ws.Range("D15").Value = 32
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = $false
$ws.Range("H15").Value = 51
$ws.Range("J15").Value = "best"
$ws.Range("K15").Value = 0.2833333333333333
$ws.Range("L15").Value = 34
$ws.Range("M15").Value = 34
$ws.Range("N15").Value = 1
$ws.Range("O15").Value = 0
$ws.Range("P15").Value = $false
$ws.Range("Q15").Value = 86

# Row 16
$ws.Range("A16").Value = "apart"
$ws.Range("B16").Value = 0.3578947368421053
$ws.Range("C16").Value = 34
$ws.Range("D16").Value = 34
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = $false
$ws.Range("H16").Value = 61
$ws.Range("J16").Value = "learn"
$ws.Range("K16").Value = 0.234375
$ws.Range("L16").Value = 30
$ws.Range("M16").Value = 30
$ws.Range("N16").Value = 1
$ws.Range("O16").Value = 0
$ws.Range("P16").Value = $false
$ws.Range("Q16").Value = 98

# Row 17
$ws.Range("A17").Value = "difficult"
$ws.Range("B17").Value = 0.3258426966292135
$ws.Range("C17").Value = 29
$ws.Range("D17").Value = 29
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = $false
$ws.Range("H17").Value = 60
$ws.Range("J17").Value = "happy"
$ws.Range("K17").Value = 0.2097902097902098
$ws.Range("L17").Value = 30
$ws.Range("M17").Value = 30
$ws.Range("N17").Value = 1
$ws.Range("O17").Value = 0
$ws.Range("P17").Value = $false
$ws.Range("Q17").Value = 113

# Row 18
$ws.Range("A18").Value = "ok"
$ws.Range("B18").Value = 0.3125
$ws.Range("C18").Value = 40
$ws.Range("D18").Value = 40
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = $false
$ws.Range("H18").Value = 88
$ws.Range("J18").Value = "enjoy"
$ws.Range("K18").Value = 0.1774193548387097
$ws.Range("L18").Value = 33
$ws.Range("M18").Value = 33
$ws.Range("N18").Value = 1
$ws.Range("O18").Value = 0
$ws.Range("P18").Value = $false
$ws.Range("Q18").Value = 153

# Row 19
$ws.Range("A19").Value = "thought"
$ws.Range("B19").Value = 0.2623762376237624
$ws.Range("C19").Value = 53
$ws.Range("D19").Value = 53
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = $false
$ws.Range("H19").Value = 149
$ws.Range("J19").Value = "christmas"
$ws.Range("K19").Value = 0.1566265060240964
$ws.Range("L19").Value = 39
$ws.Range("M19").Value = 39
$ws.Range("N19").Value = 1
$ws.Range("O19").Value = 0
$ws.Range("P19").Value = $false
$ws.Range("Q19").Value = 210

# Row 20
$ws.Range("A20").Value = "size"
$ws.Range("B20").Value = 0.2216494845360825
$ws.Range("C20").Value = 43
$ws.Range("D20").Value = 43
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = $false
$ws.Range("H20").Value = 151
$ws.Range("J20").Value = "fun"
$ws.Range("K20").Value = 0.1454864154250657
$ws.Range("L20").Value = 166
$ws.Range("M20").Value = 166
$ws.Range("N20").Value = 1
$ws.Range("O20").Value = 0
$ws.Range("P20").Value = $false
$ws.Range("Q20").Value = 975

# Row 21
$ws.Range("A21").Value = "item"
$ws.Range("B21").Value = 0.1739130434782609
$ws.Range("C21").Value = 48
$ws.Range("D21").Value = 48
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = $false
$ws.Range("H21").Value = 228
$ws.Range("J21").Value = "game"
$ws.Range("K21").Value = 0.07082521117608837
$ws.Range("L21").Value = 109
$ws.Range("M21").Value = 111
$ws.Range("N21").Value = 0.98
$ws.Range("O21").Value = 0.02000000000000002
$ws.Range("P21").Value = $true
$ws.Range("Q21").Value = 1430

# --- Row 22: update A-H, and ADD new J-Q block (copy style from J21:Q21 first) ---
$ws.Range("J21:Q21").Copy($ws.Range("J22:Q22"))
$ws.Range("A22").Value = "hard"
$ws.Range("B22").Value = 0.17
$ws.Range("C22").Value = 34
$ws.Range("D22").Value = 34
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = $false
$ws.Range("H22").Value = 166
$ws.Range("J22").Value = "play"
$ws.Range("K22").Value = 0.044
$ws.Range("L22").Value = 33
$ws.Range("M22").Value = 35
$ws.Range("N22").Value = 0.94
$ws.Range("O22").Value = 0.06000000000000005
$ws.Range("P22").Value = $true
$ws.Range("Q22").Value = 717

# --- Rows 23-33: update A-H only ---
# Row 23
$ws.Range("A23").Value = "money"
$ws.Range("B23").Value = 0.1645569620253164
$ws.Range("C23").Value = 52
$ws.Range("D23").Value = 52
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 1
$ws.Range("G23").Value = $false
$ws.Range("H23").Value = 264

# Row 24
$ws.Range("A24").Value = "would"
$ws.Range("B24").Value = 0.150521609538003
$ws.Range("C24").Value = 101
$ws.Range("D24").Value = 104
$ws.Range("E24").Value = 0.03
$ws.Range("F24").Value = 0.97
$ws.Range("G24").Value = $true
$ws.Range("H24").Value = 570

# Row 25
$ws.Range("A25").Value = "price"
$ws.Range("B25").Value = 0.1498559077809798
$ws.Range("C25").Value = 52
$ws.Range("D25").Value = 53
$ws.Range("E25").Value = 0.02
$ws.Range("F25").Value = 0.98
$ws.Range("G25").Value = $true
$ws.Range("H25").Value = 295

# Row 26
$ws.Range("A26").Value = "better"
$ws.Range("B26").Value = 0.1495327102803738
$ws.Range("C26").Value = 32
$ws.Range("D26").Value = 32
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 1
$ws.Range("G26").Value = $false
$ws.Range("H26").Value = 182

# Row 27
$ws.Range("A27").Value = "work"
$ws.Range("B27").Value = 0.1487341772151899
$ws.Range("C27").Value = 47
$ws.Range("D27").Value = 47
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = $false
$ws.Range("H27").Value = 269

# Row 28
Set-TextValue $ws.Range("A28") "3"
$ws.Range("B28").Value = 0.1209677419354839
$ws.Range("C28").Value = 30
$ws.Range("D28").Value = 30
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 1
$ws.Range("G28").Value = $false
$ws.Range("H28").Value = 218

# Row 29
Set-TextValue $ws.Range("A29") "2"
$ws.Range("B29").Value = 0.1198501872659176
$ws.Range("C29").Value = 32
$ws.Range("D29").Value = 32
$ws.Range("E29").Value = 0
$ws.Range("F29").Value = 1
$ws.Range("G29").Value = $false
$ws.Range("H29").Value = 235

# Row 30
$ws.Range("A30").Value = "product"
$ws.Range("B30").Value = 0.1167400881057269
$ws.Range("C30").Value = 53
$ws.Range("D30").Value = 53
$ws.Range("E30").Value = 0
$ws.Range("F30").Value = 1
$ws.Range("G30").Value = $false
$ws.Range("H30").Value = 401

# Row 31
$ws.Range("A31").Value = "buy"
$ws.Range("B31").Value = 0.08169014084507042
$ws.Range("C31").Value = 29
$ws.Range("D31").Value = 29
$ws.Range("E31").Value = 0
$ws.Range("F31").Value = 1
$ws.Range("G31").Value = $false
$ws.Range("H31").Value = 326

# Row 32
$ws.Range("A32").Value = "little"
$ws.Range("B32").Value = 0.0757238307349666
$ws.Range("C32").Value = 34
$ws.Range("D32").Value = 34
$ws.Range("E32").Value = 0
$ws.Range("F32").Value = 1
$ws.Range("G32").Value = $false
$ws.Range("H32").Value = 415

# Row 33
$ws.Range("A33").Value = "like"
$ws.Range("B33").Value = 0.06765676567656766
$ws.Range("C33").Value = 41
$ws.Range("D33").Value = 43
$ws.Range("E33").Value = 0.05
$ws.Range("F33").Value = 0.95
$ws.Range("G33").Value = $true
$ws.Range("H33").Value = 565

# --- Row 34: brand-new row (copy style/border from A33 first) ---
$ws.Range("A33").Copy($ws.Range("A34"))
$ws.Range("A34").Value = "one"
$ws.Range("B34").Value = 0.04336734693877551
$ws.Range("C34").Value = 34
$ws.Range("D34").Value = 44
$ws.Range("E34").Value = 0.23
$ws.Range("F34").Value = 0.77
$ws.Range("G34").Value = $true
$ws.Range("H34").Value = 750

Write-Host "edit complete"